# Applies the weekly Femacal de La Calera / Zapallo italiano price update:
# two new rows (305, 306) are inserted with fresh readings, and the rows that
# previously occupied 305-333 shift down to 307-335. Because we know the exact
# target content for every row in the affected range, the simplest and most
# robust way to reach that state through the Excel object model is to write the
# final values directly into each cell (rather than performing a literal
# Rows.Insert shift), and to restore the date-column (D) number format on the
# two rows that are brand new to the sheet (305-335 overall; 334-335 previously
# did not exist in the workbook at all).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 305
$ws.Cells.Item(305, 1).Value = 3
$ws.Cells.Item(305, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(305, 3).Value = "Coquimbo"
$ws.Cells.Item(305, 4).Value = 44578
$ws.Cells.Item(305, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(305, 5).Value = 5
$ws.Cells.Item(305, 6).Value = 100112032
$ws.Cells.Item(305, 7).Value = "Zapallo italiano"
$ws.Cells.Item(305, 8).Value = "Sin especificar"
$ws.Cells.Item(305, 9).Value = "Primera"
$ws.Cells.Item(305, 10).Value = 140
$ws.Cells.Item(305, 11).Value = 4500
$ws.Cells.Item(305, 12).Value = 5000
$ws.Cells.Item(305, 13).Value = 4714
$ws.Cells.Item(305, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(305, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(305, 16).Value = 131
$ws.Cells.Item(305, 17).Value = 36
$ws.Cells.Item(305, 18).Value = "Hortaliza"

# Row 306
$ws.Cells.Item(306, 1).Value = 3
$ws.Cells.Item(306, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(306, 3).Value = "Coquimbo"
$ws.Cells.Item(306, 4).Value = 44578
$ws.Cells.Item(306, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(306, 5).Value = 5
$ws.Cells.Item(306, 6).Value = 100112032
$ws.Cells.Item(306, 7).Value = "Zapallo italiano"
$ws.Cells.Item(306, 8).Value = "Sin especificar"
$ws.Cells.Item(306, 9).Value = "Primera"
$ws.Cells.Item(306, 10).Value = 200
$ws.Cells.Item(306, 11).Value = 8000
$ws.Cells.Item(306, 12).Value = 8500
$ws.Cells.Item(306, 13).Value = 8175
$ws.Cells.Item(306, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(306, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(306, 16).Value = 117
$ws.Cells.Item(306, 17).Value = 70
$ws.Cells.Item(306, 18).Value = "Hortaliza"

# Row 307
$ws.Cells.Item(307, 1).Value = 3
$ws.Cells.Item(307, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(307, 3).Value = "Coquimbo"
$ws.Cells.Item(307, 4).Value = 44490
$ws.Cells.Item(307, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(307, 5).Value = 5
$ws.Cells.Item(307, 6).Value = 100112032
$ws.Cells.Item(307, 7).Value = "Zapallo italiano"
$ws.Cells.Item(307, 8).Value = "Sin especificar"
$ws.Cells.Item(307, 9).Value = "Primera"
$ws.Cells.Item(307, 10).Value = 105
$ws.Cells.Item(307, 11).Value = 8000
$ws.Cells.Item(307, 12).Value = 8500
$ws.Cells.Item(307, 13).Value = 8262
$ws.Cells.Item(307, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(307, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(307, 16).Value = 230
$ws.Cells.Item(307, 17).Value = 36
$ws.Cells.Item(307, 18).Value = "Hortaliza"

# Row 308
$ws.Cells.Item(308, 1).Value = 3
$ws.Cells.Item(308, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(308, 3).Value = "Coquimbo"
$ws.Cells.Item(308, 4).Value = 44490
$ws.Cells.Item(308, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(308, 5).Value = 5
$ws.Cells.Item(308, 6).Value = 100112032
$ws.Cells.Item(308, 7).Value = "Zapallo italiano"
$ws.Cells.Item(308, 8).Value = "Sin especificar"
$ws.Cells.Item(308, 9).Value = "Primera"
$ws.Cells.Item(308, 10).Value = 100
$ws.Cells.Item(308, 11).Value = 9500
$ws.Cells.Item(308, 12).Value = 10000
$ws.Cells.Item(308, 13).Value = 9750
$ws.Cells.Item(308, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(308, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(308, 16).Value = 139
$ws.Cells.Item(308, 17).Value = 70
$ws.Cells.Item(308, 18).Value = "Hortaliza"

# Row 309
$ws.Cells.Item(309, 1).Value = 3
$ws.Cells.Item(309, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(309, 3).Value = "Coquimbo"
$ws.Cells.Item(309, 4).Value = 44427
$ws.Cells.Item(309, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(309, 5).Value = 5
$ws.Cells.Item(309, 6).Value = 100112032
$ws.Cells.Item(309, 7).Value = "Zapallo italiano"
$ws.Cells.Item(309, 8).Value = "Sin especificar"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 225
$ws.Cells.Item(309, 11).Value = 7000
$ws.Cells.Item(309, 12).Value = 8000
$ws.Cells.Item(309, 13).Value = 7522
$ws.Cells.Item(309, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(309, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(309, 16).Value = 107
$ws.Cells.Item(309, 17).Value = 70
$ws.Cells.Item(309, 18).Value = "Hortaliza"

# Row 310
$ws.Cells.Item(310, 1).Value = 3
$ws.Cells.Item(310, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(310, 3).Value = "Coquimbo"
$ws.Cells.Item(310, 4).Value = 44491
$ws.Cells.Item(310, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(310, 5).Value = 5
$ws.Cells.Item(310, 6).Value = 100112032
$ws.Cells.Item(310, 7).Value = "Zapallo italiano"
$ws.Cells.Item(310, 8).Value = "Sin especificar"
$ws.Cells.Item(310, 9).Value = "Primera"
$ws.Cells.Item(310, 10).Value = 90
$ws.Cells.Item(310, 11).Value = 9500
$ws.Cells.Item(310, 12).Value = 10000
$ws.Cells.Item(310, 13).Value = 9750
$ws.Cells.Item(310, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(310, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(310, 16).Value = 139
$ws.Cells.Item(310, 17).Value = 70
$ws.Cells.Item(310, 18).Value = "Hortaliza"

# Row 311
$ws.Cells.Item(311, 1).Value = 3
$ws.Cells.Item(311, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(311, 3).Value = "Coquimbo"
$ws.Cells.Item(311, 4).Value = 44293
$ws.Cells.Item(311, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(311, 5).Value = 5
$ws.Cells.Item(311, 6).Value = 100112032
$ws.Cells.Item(311, 7).Value = "Zapallo italiano"
$ws.Cells.Item(311, 8).Value = "Sin especificar"
$ws.Cells.Item(311, 9).Value = "Primera"
$ws.Cells.Item(311, 10).Value = 190
$ws.Cells.Item(311, 11).Value = 8000
$ws.Cells.Item(311, 12).Value = 8500
$ws.Cells.Item(311, 13).Value = 8158
$ws.Cells.Item(311, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(311, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(311, 16).Value = 117
$ws.Cells.Item(311, 17).Value = 70
$ws.Cells.Item(311, 18).Value = "Hortaliza"

# Row 312
$ws.Cells.Item(312, 1).Value = 3
$ws.Cells.Item(312, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(312, 3).Value = "Coquimbo"
$ws.Cells.Item(312, 4).Value = 44266
$ws.Cells.Item(312, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(312, 5).Value = 5
$ws.Cells.Item(312, 6).Value = 100112032
$ws.Cells.Item(312, 7).Value = "Zapallo italiano"
$ws.Cells.Item(312, 8).Value = "Sin especificar"
$ws.Cells.Item(312, 9).Value = "Primera"
$ws.Cells.Item(312, 10).Value = 160
$ws.Cells.Item(312, 11).Value = 8000
$ws.Cells.Item(312, 12).Value = 9000
$ws.Cells.Item(312, 13).Value = 8500
$ws.Cells.Item(312, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(312, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(312, 16).Value = 121
$ws.Cells.Item(312, 17).Value = 70
$ws.Cells.Item(312, 18).Value = "Hortaliza"

# Row 313
$ws.Cells.Item(313, 1).Value = 3
$ws.Cells.Item(313, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(313, 3).Value = "Coquimbo"
$ws.Cells.Item(313, 4).Value = 44533
$ws.Cells.Item(313, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(313, 5).Value = 5
$ws.Cells.Item(313, 6).Value = 100112032
$ws.Cells.Item(313, 7).Value = "Zapallo italiano"
$ws.Cells.Item(313, 8).Value = "Sin especificar"
$ws.Cells.Item(313, 9).Value = "Primera"
$ws.Cells.Item(313, 10).Value = 120
$ws.Cells.Item(313, 11).Value = 4500
$ws.Cells.Item(313, 12).Value = 5000
$ws.Cells.Item(313, 13).Value = 4750
$ws.Cells.Item(313, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(313, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(313, 16).Value = 132
$ws.Cells.Item(313, 17).Value = 36
$ws.Cells.Item(313, 18).Value = "Hortaliza"

# Row 314
$ws.Cells.Item(314, 1).Value = 3
$ws.Cells.Item(314, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(314, 3).Value = "Coquimbo"
$ws.Cells.Item(314, 4).Value = 44533
$ws.Cells.Item(314, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(314, 5).Value = 5
$ws.Cells.Item(314, 6).Value = 100112032
$ws.Cells.Item(314, 7).Value = "Zapallo italiano"
$ws.Cells.Item(314, 8).Value = "Sin especificar"
$ws.Cells.Item(314, 9).Value = "Primera"
$ws.Cells.Item(314, 10).Value = 185
$ws.Cells.Item(314, 11).Value = 7500
$ws.Cells.Item(314, 12).Value = 8000
$ws.Cells.Item(314, 13).Value = 7768
$ws.Cells.Item(314, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(314, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(314, 16).Value = 111
$ws.Cells.Item(314, 17).Value = 70
$ws.Cells.Item(314, 18).Value = "Hortaliza"

# Row 315
$ws.Cells.Item(315, 1).Value = 3
$ws.Cells.Item(315, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(315, 3).Value = "Coquimbo"
$ws.Cells.Item(315, 4).Value = 44264
$ws.Cells.Item(315, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(315, 5).Value = 5
$ws.Cells.Item(315, 6).Value = 100112032
$ws.Cells.Item(315, 7).Value = "Zapallo italiano"
$ws.Cells.Item(315, 8).Value = "Sin especificar"
$ws.Cells.Item(315, 9).Value = "Primera"
$ws.Cells.Item(315, 10).Value = 114
$ws.Cells.Item(315, 11).Value = 9000
$ws.Cells.Item(315, 12).Value = 9000
$ws.Cells.Item(315, 13).Value = 9000
$ws.Cells.Item(315, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(315, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(315, 16).Value = 129
$ws.Cells.Item(315, 17).Value = 70
$ws.Cells.Item(315, 18).Value = "Hortaliza"

# Row 316
$ws.Cells.Item(316, 1).Value = 3
$ws.Cells.Item(316, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(316, 3).Value = "Coquimbo"
$ws.Cells.Item(316, 4).Value = 44494
$ws.Cells.Item(316, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(316, 5).Value = 5
$ws.Cells.Item(316, 6).Value = 100112032
$ws.Cells.Item(316, 7).Value = "Zapallo italiano"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 130
$ws.Cells.Item(316, 11).Value = 8000
$ws.Cells.Item(316, 12).Value = 8500
$ws.Cells.Item(316, 13).Value = 8269
$ws.Cells.Item(316, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(316, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(316, 16).Value = 230
$ws.Cells.Item(316, 17).Value = 36
$ws.Cells.Item(316, 18).Value = "Hortaliza"

# Row 317
$ws.Cells.Item(317, 1).Value = 3
$ws.Cells.Item(317, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(317, 3).Value = "Coquimbo"
$ws.Cells.Item(317, 4).Value = 44494
$ws.Cells.Item(317, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(317, 5).Value = 5
$ws.Cells.Item(317, 6).Value = 100112032
$ws.Cells.Item(317, 7).Value = "Zapallo italiano"
$ws.Cells.Item(317, 8).Value = "Sin especificar"
$ws.Cells.Item(317, 9).Value = "Primera"
$ws.Cells.Item(317, 10).Value = 130
$ws.Cells.Item(317, 11).Value = 9000
$ws.Cells.Item(317, 12).Value = 9500
$ws.Cells.Item(317, 13).Value = 9231
$ws.Cells.Item(317, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(317, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(317, 16).Value = 132
$ws.Cells.Item(317, 17).Value = 70
$ws.Cells.Item(317, 18).Value = "Hortaliza"

# Row 318
$ws.Cells.Item(318, 1).Value = 3
$ws.Cells.Item(318, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(318, 3).Value = "Coquimbo"
$ws.Cells.Item(318, 4).Value = 44571
$ws.Cells.Item(318, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(318, 5).Value = 5
$ws.Cells.Item(318, 6).Value = 100112032
$ws.Cells.Item(318, 7).Value = "Zapallo italiano"
$ws.Cells.Item(318, 8).Value = "Sin especificar"
$ws.Cells.Item(318, 9).Value = "Primera"
$ws.Cells.Item(318, 10).Value = 190
$ws.Cells.Item(318, 11).Value = 4500
$ws.Cells.Item(318, 12).Value = 5000
$ws.Cells.Item(318, 13).Value = 4763
$ws.Cells.Item(318, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(318, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(318, 16).Value = 132
$ws.Cells.Item(318, 17).Value = 36
$ws.Cells.Item(318, 18).Value = "Hortaliza"

# Row 319
$ws.Cells.Item(319, 1).Value = 3
$ws.Cells.Item(319, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(319, 3).Value = "Coquimbo"
$ws.Cells.Item(319, 4).Value = 44571
$ws.Cells.Item(319, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(319, 5).Value = 5
$ws.Cells.Item(319, 6).Value = 100112032
$ws.Cells.Item(319, 7).Value = "Zapallo italiano"
$ws.Cells.Item(319, 8).Value = "Sin especificar"
$ws.Cells.Item(319, 9).Value = "Primera"
$ws.Cells.Item(319, 10).Value = 130
$ws.Cells.Item(319, 11).Value = 8000
$ws.Cells.Item(319, 12).Value = 8500
$ws.Cells.Item(319, 13).Value = 8231
$ws.Cells.Item(319, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(319, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(319, 16).Value = 118
$ws.Cells.Item(319, 17).Value = 70
$ws.Cells.Item(319, 18).Value = "Hortaliza"

# Row 320
$ws.Cells.Item(320, 1).Value = 3
$ws.Cells.Item(320, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(320, 3).Value = "Coquimbo"
$ws.Cells.Item(320, 4).Value = 44390
$ws.Cells.Item(320, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(320, 5).Value = 5
$ws.Cells.Item(320, 6).Value = 100112032
$ws.Cells.Item(320, 7).Value = "Zapallo italiano"
$ws.Cells.Item(320, 8).Value = "Sin especificar"
$ws.Cells.Item(320, 9).Value = "Primera"
$ws.Cells.Item(320, 10).Value = 45
$ws.Cells.Item(320, 11).Value = 12000
$ws.Cells.Item(320, 12).Value = 12000
$ws.Cells.Item(320, 13).Value = 12000
$ws.Cells.Item(320, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(320, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(320, 16).Value = 171
$ws.Cells.Item(320, 17).Value = 70
$ws.Cells.Item(320, 18).Value = "Hortaliza"

# Row 321
$ws.Cells.Item(321, 1).Value = 3
$ws.Cells.Item(321, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(321, 3).Value = "Coquimbo"
$ws.Cells.Item(321, 4).Value = 44279
$ws.Cells.Item(321, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(321, 5).Value = 5
$ws.Cells.Item(321, 6).Value = 100112032
$ws.Cells.Item(321, 7).Value = "Zapallo italiano"
$ws.Cells.Item(321, 8).Value = "Sin especificar"
$ws.Cells.Item(321, 9).Value = "Primera"
$ws.Cells.Item(321, 10).Value = 170
$ws.Cells.Item(321, 11).Value = 7000
$ws.Cells.Item(321, 12).Value = 9000
$ws.Cells.Item(321, 13).Value = 8147
$ws.Cells.Item(321, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(321, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(321, 16).Value = 116
$ws.Cells.Item(321, 17).Value = 70
$ws.Cells.Item(321, 18).Value = "Hortaliza"

# Row 322
$ws.Cells.Item(322, 1).Value = 3
$ws.Cells.Item(322, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(322, 3).Value = "Coquimbo"
$ws.Cells.Item(322, 4).Value = 44481
$ws.Cells.Item(322, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(322, 5).Value = 5
$ws.Cells.Item(322, 6).Value = 100112032
$ws.Cells.Item(322, 7).Value = "Zapallo italiano"
$ws.Cells.Item(322, 8).Value = "Sin especificar"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 100
$ws.Cells.Item(322, 11).Value = 13000
$ws.Cells.Item(322, 12).Value = 14000
$ws.Cells.Item(322, 13).Value = 13500
$ws.Cells.Item(322, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(322, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(322, 16).Value = 193
$ws.Cells.Item(322, 17).Value = 70
$ws.Cells.Item(322, 18).Value = "Hortaliza"

# Row 323
$ws.Cells.Item(323, 1).Value = 3
$ws.Cells.Item(323, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(323, 3).Value = "Coquimbo"
$ws.Cells.Item(323, 4).Value = 44277
$ws.Cells.Item(323, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(323, 5).Value = 5
$ws.Cells.Item(323, 6).Value = 100112032
$ws.Cells.Item(323, 7).Value = "Zapallo italiano"
$ws.Cells.Item(323, 8).Value = "Sin especificar"
$ws.Cells.Item(323, 9).Value = "Primera"
$ws.Cells.Item(323, 10).Value = 140
$ws.Cells.Item(323, 11).Value = 9000
$ws.Cells.Item(323, 12).Value = 9000
$ws.Cells.Item(323, 13).Value = 9000
$ws.Cells.Item(323, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(323, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(323, 16).Value = 129
$ws.Cells.Item(323, 17).Value = 70
$ws.Cells.Item(323, 18).Value = "Hortaliza"

# Row 324
$ws.Cells.Item(324, 1).Value = 3
$ws.Cells.Item(324, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(324, 3).Value = "Coquimbo"
$ws.Cells.Item(324, 4).Value = 44525
$ws.Cells.Item(324, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(324, 5).Value = 5
$ws.Cells.Item(324, 6).Value = 100112032
$ws.Cells.Item(324, 7).Value = "Zapallo italiano"
$ws.Cells.Item(324, 8).Value = "Sin especificar"
$ws.Cells.Item(324, 9).Value = "Primera"
$ws.Cells.Item(324, 10).Value = 160
$ws.Cells.Item(324, 11).Value = 4500
$ws.Cells.Item(324, 12).Value = 4800
$ws.Cells.Item(324, 13).Value = 4650
$ws.Cells.Item(324, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(324, 15).Value = "Limache"
$ws.Cells.Item(324, 16).Value = 129
$ws.Cells.Item(324, 17).Value = 36
$ws.Cells.Item(324, 18).Value = "Hortaliza"

# Row 325
$ws.Cells.Item(325, 1).Value = 3
$ws.Cells.Item(325, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(325, 3).Value = "Coquimbo"
$ws.Cells.Item(325, 4).Value = 44525
$ws.Cells.Item(325, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(325, 5).Value = 5
$ws.Cells.Item(325, 6).Value = 100112032
$ws.Cells.Item(325, 7).Value = "Zapallo italiano"
$ws.Cells.Item(325, 8).Value = "Sin especificar"
$ws.Cells.Item(325, 9).Value = "Primera"
$ws.Cells.Item(325, 10).Value = 245
$ws.Cells.Item(325, 11).Value = 7000
$ws.Cells.Item(325, 12).Value = 7500
$ws.Cells.Item(325, 13).Value = 7261
$ws.Cells.Item(325, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(325, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(325, 16).Value = 104
$ws.Cells.Item(325, 17).Value = 70
$ws.Cells.Item(325, 18).Value = "Hortaliza"

# Row 326
$ws.Cells.Item(326, 1).Value = 3
$ws.Cells.Item(326, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(326, 3).Value = "Coquimbo"
$ws.Cells.Item(326, 4).Value = 44327
$ws.Cells.Item(326, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(326, 5).Value = 5
$ws.Cells.Item(326, 6).Value = 100112032
$ws.Cells.Item(326, 7).Value = "Zapallo italiano"
$ws.Cells.Item(326, 8).Value = "Sin especificar"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 220
$ws.Cells.Item(326, 11).Value = 6000
$ws.Cells.Item(326, 12).Value = 6500
$ws.Cells.Item(326, 13).Value = 6255
$ws.Cells.Item(326, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(326, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(326, 16).Value = 89
$ws.Cells.Item(326, 17).Value = 70
$ws.Cells.Item(326, 18).Value = "Hortaliza"

# Row 327
$ws.Cells.Item(327, 1).Value = 3
$ws.Cells.Item(327, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(327, 3).Value = "Coquimbo"
$ws.Cells.Item(327, 4).Value = 44503
$ws.Cells.Item(327, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(327, 5).Value = 5
$ws.Cells.Item(327, 6).Value = 100112032
$ws.Cells.Item(327, 7).Value = "Zapallo italiano"
$ws.Cells.Item(327, 8).Value = "Sin especificar"
$ws.Cells.Item(327, 9).Value = "Primera"
$ws.Cells.Item(327, 10).Value = 75
$ws.Cells.Item(327, 11).Value = 7000
$ws.Cells.Item(327, 12).Value = 7500
$ws.Cells.Item(327, 13).Value = 7167
$ws.Cells.Item(327, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(327, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(327, 16).Value = 199
$ws.Cells.Item(327, 17).Value = 36
$ws.Cells.Item(327, 18).Value = "Hortaliza"

# Row 328
$ws.Cells.Item(328, 1).Value = 3
$ws.Cells.Item(328, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(328, 3).Value = "Coquimbo"
$ws.Cells.Item(328, 4).Value = 44503
$ws.Cells.Item(328, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(328, 5).Value = 5
$ws.Cells.Item(328, 6).Value = 100112032
$ws.Cells.Item(328, 7).Value = "Zapallo italiano"
$ws.Cells.Item(328, 8).Value = "Sin especificar"
$ws.Cells.Item(328, 9).Value = "Primera"
$ws.Cells.Item(328, 10).Value = 50
$ws.Cells.Item(328, 11).Value = 8500
$ws.Cells.Item(328, 12).Value = 8500
$ws.Cells.Item(328, 13).Value = 8500
$ws.Cells.Item(328, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(328, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(328, 16).Value = 121
$ws.Cells.Item(328, 17).Value = 70
$ws.Cells.Item(328, 18).Value = "Hortaliza"

# Row 329
$ws.Cells.Item(329, 1).Value = 3
$ws.Cells.Item(329, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(329, 3).Value = "Coquimbo"
$ws.Cells.Item(329, 4).Value = 44462
$ws.Cells.Item(329, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(329, 5).Value = 5
$ws.Cells.Item(329, 6).Value = 100112032
$ws.Cells.Item(329, 7).Value = "Zapallo italiano"
$ws.Cells.Item(329, 8).Value = "Sin especificar"
$ws.Cells.Item(329, 9).Value = "Primera"
$ws.Cells.Item(329, 10).Value = 205
$ws.Cells.Item(329, 11).Value = 11000
$ws.Cells.Item(329, 12).Value = 12000
$ws.Cells.Item(329, 13).Value = 11512
$ws.Cells.Item(329, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(329, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(329, 16).Value = 164
$ws.Cells.Item(329, 17).Value = 70
$ws.Cells.Item(329, 18).Value = "Hortaliza"

# Row 330
$ws.Cells.Item(330, 1).Value = 3
$ws.Cells.Item(330, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(330, 3).Value = "Coquimbo"
$ws.Cells.Item(330, 4).Value = 44384
$ws.Cells.Item(330, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(330, 5).Value = 5
$ws.Cells.Item(330, 6).Value = 100112032
$ws.Cells.Item(330, 7).Value = "Zapallo italiano"
$ws.Cells.Item(330, 8).Value = "Sin especificar"
$ws.Cells.Item(330, 9).Value = "Primera"
$ws.Cells.Item(330, 10).Value = 115
$ws.Cells.Item(330, 11).Value = 9000
$ws.Cells.Item(330, 12).Value = 9500
$ws.Cells.Item(330, 13).Value = 9261
$ws.Cells.Item(330, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(330, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(330, 16).Value = 132
$ws.Cells.Item(330, 17).Value = 70
$ws.Cells.Item(330, 18).Value = "Hortaliza"

# Row 331
$ws.Cells.Item(331, 1).Value = 3
$ws.Cells.Item(331, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(331, 3).Value = "Coquimbo"
$ws.Cells.Item(331, 4).Value = 44512
$ws.Cells.Item(331, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(331, 5).Value = 5
$ws.Cells.Item(331, 6).Value = 100112032
$ws.Cells.Item(331, 7).Value = "Zapallo italiano"
$ws.Cells.Item(331, 8).Value = "Sin especificar"
$ws.Cells.Item(331, 9).Value = "Primera"
$ws.Cells.Item(331, 10).Value = 125
$ws.Cells.Item(331, 11).Value = 4000
$ws.Cells.Item(331, 12).Value = 4500
$ws.Cells.Item(331, 13).Value = 4240
$ws.Cells.Item(331, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(331, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(331, 16).Value = 118
$ws.Cells.Item(331, 17).Value = 36
$ws.Cells.Item(331, 18).Value = "Hortaliza"

# Row 332
$ws.Cells.Item(332, 1).Value = 3
$ws.Cells.Item(332, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(332, 3).Value = "Coquimbo"
$ws.Cells.Item(332, 4).Value = 44512
$ws.Cells.Item(332, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(332, 5).Value = 5
$ws.Cells.Item(332, 6).Value = 100112032
$ws.Cells.Item(332, 7).Value = "Zapallo italiano"
$ws.Cells.Item(332, 8).Value = "Sin especificar"
$ws.Cells.Item(332, 9).Value = "Primera"
$ws.Cells.Item(332, 10).Value = 128
$ws.Cells.Item(332, 11).Value = 7000
$ws.Cells.Item(332, 12).Value = 7500
$ws.Cells.Item(332, 13).Value = 7234
$ws.Cells.Item(332, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(332, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(332, 16).Value = 103
$ws.Cells.Item(332, 17).Value = 70
$ws.Cells.Item(332, 18).Value = "Hortaliza"

# Row 333
$ws.Cells.Item(333, 1).Value = 3
$ws.Cells.Item(333, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(333, 3).Value = "Coquimbo"
$ws.Cells.Item(333, 4).Value = 44312
$ws.Cells.Item(333, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(333, 5).Value = 5
$ws.Cells.Item(333, 6).Value = 100112032
$ws.Cells.Item(333, 7).Value = "Zapallo italiano"
$ws.Cells.Item(333, 8).Value = "Sin especificar"
$ws.Cells.Item(333, 9).Value = "Primera"
$ws.Cells.Item(333, 10).Value = 105
$ws.Cells.Item(333, 11).Value = 6000
$ws.Cells.Item(333, 12).Value = 6500
$ws.Cells.Item(333, 13).Value = 6238
$ws.Cells.Item(333, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(333, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(333, 16).Value = 89
$ws.Cells.Item(333, 17).Value = 70
$ws.Cells.Item(333, 18).Value = "Hortaliza"

# Row 334
$ws.Cells.Item(334, 1).Value = 3
$ws.Cells.Item(334, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(334, 3).Value = "Coquimbo"
$ws.Cells.Item(334, 4).Value = 44511
$ws.Cells.Item(334, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(334, 5).Value = 5
$ws.Cells.Item(334, 6).Value = 100112032
$ws.Cells.Item(334, 7).Value = "Zapallo italiano"
$ws.Cells.Item(334, 8).Value = "Sin especificar"
$ws.Cells.Item(334, 9).Value = "Primera"
$ws.Cells.Item(334, 10).Value = 120
$ws.Cells.Item(334, 11).Value = 4000
$ws.Cells.Item(334, 12).Value = 4500
$ws.Cells.Item(334, 13).Value = 4250
$ws.Cells.Item(334, 14).Value = "`$/caja 36 unidades"
$ws.Cells.Item(334, 15).Value = "Limache"
$ws.Cells.Item(334, 16).Value = 118
$ws.Cells.Item(334, 17).Value = 36
$ws.Cells.Item(334, 18).Value = "Hortaliza"

# Row 335
$ws.Cells.Item(335, 1).Value = 3
$ws.Cells.Item(335, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(335, 3).Value = "Coquimbo"
$ws.Cells.Item(335, 4).Value = 44511
$ws.Cells.Item(335, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(335, 5).Value = 5
$ws.Cells.Item(335, 6).Value = 100112032
$ws.Cells.Item(335, 7).Value = "Zapallo italiano"
$ws.Cells.Item(335, 8).Value = "Sin especificar"
$ws.Cells.Item(335, 9).Value = "Primera"
$ws.Cells.Item(335, 10).Value = 230
$ws.Cells.Item(335, 11).Value = 7000
$ws.Cells.Item(335, 12).Value = 7500
$ws.Cells.Item(335, 13).Value = 7261
$ws.Cells.Item(335, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(335, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(335, 16).Value = 104
$ws.Cells.Item(335, 17).Value = 70
$ws.Cells.Item(335, 18).Value = "Hortaliza"

